# Apply "Generate Report for Handback" timestamp updates across the
# Overview, zh-cn, and de-de worksheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for the first file row.
$wsOverview.Range("G2").Value = "2016-10-21 00:57:54"

# zh-cn sheet: Correspond Handoff Datetime (H2) and Correspond Handback DateTime (K2)
$wsZhCn.Range("H2").Value = "2016-10-21 00:57:43"
$wsZhCn.Range("K2").Value = "2016-10-21 00:58:23"

# de-de sheet: Correspond Handoff Datetime (H2) matches the Overview's updated date,
# Correspond Handback DateTime (K2) gets its own new value.
$wsDeDe.Range("H2").Value = "2016-10-21 00:57:54"
$wsDeDe.Range("K2").Value = "2016-10-21 00:58:41"
